$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 114
$ws.Range("B1").Value = 160.3999999999996
$ws.Range("C1").Value = 293.8000140938362
$ws.Range("A2").Value = 114
